$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 6577.2
$ws.Range("J51").Value = 3000
$ws.Range("L51").Value = 3000
$ws.Range("N51").Value = -3968

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 8999.833000000001
$ws.Range("I70").Value = 5000
$ws.Range("J70").Value = 9799.799999999999
$ws.Range("K70").Value = 15000
$ws.Range("L70").Value = 29399.4
$ws.Range("M70").Value = -14730
$ws.Range("N70").Value = -29939.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 8999.833000000001
$ws.Range("I73").Value = 5000
$ws.Range("J73").Value = 9799.799999999999
$ws.Range("K73").Value = 15000
$ws.Range("L73").Value = 29399.4
$ws.Range("M73").Value = -14064
$ws.Range("N73").Value = -31271.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1479.0526
$ws.Range("J129").Value = 2517.375
$ws.Range("L129").Value = 7552.125
$ws.Range("N129").Value = -17552.125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2412.625
$ws.Range("J138").Value = 3858.8333
$ws.Range("L138").Value = 11576.4999
$ws.Range("N138").Value = -21856.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2599
$ws.Range("I132").Value = 2498.75
$ws.Range("K132").Value = 7496.25
$ws.Range("M132").Value = -4966.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1541.625
$ws.Range("I105").Value = 1476.1428
$ws.Range("K105").Value = 1476.1428
$ws.Range("M105").Value = 270.8571999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1500.2142
$ws.Range("I22").Value = 944.8889
$ws.Range("K22").Value = 944.8889
$ws.Range("M22").Value = -594.8889

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H37").Value = 20000
$ws.Range("J37").Value = 20000
$ws.Range("L37").Value = 20000
$ws.Range("N37").Value = -20214

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 30500
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 30500
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 30500
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -31312

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H91").Value = 30500
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 30500
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 30500
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -33308

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H111").Value = 54057.332
$ws.Range("J111").Value = 54057.332
$ws.Range("L111").Value = 54057.332
$ws.Range("N111").Value = -62237.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H121").Value = 40326
$ws.Range("J121").Value = 40326
$ws.Range("L121").Value = 40326
$ws.Range("N121").Value = -42946

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3199.7368
$ws.Range("I132").Value = 2564.125
$ws.Range("J132").Value = 6589.6665
$ws.Range("K132").Value = 7692.375
$ws.Range("L132").Value = 19768.9995
$ws.Range("M132").Value = -5162.375
$ws.Range("N132").Value = -24828.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 348.63635
$ws.Range("I8").Value = 348.63635
$ws.Range("K8").Value = 1045.90905
$ws.Range("M8").Value = -906.90905

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2309.5833
$ws.Range("J34").Value = 2428.4546
$ws.Range("L34").Value = 7285.3638
$ws.Range("N34").Value = -7453.3638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 2500
$ws.Range("J42").Value = 2500
$ws.Range("L42").Value = 7500
$ws.Range("N42").Value = -8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4675.36
$ws.Range("I80").Value = 4677.55
$ws.Range("J80").Value = 4666.6
$ws.Range("K80").Value = 14032.65
$ws.Range("L80").Value = 13999.8
$ws.Range("M80").Value = -13096.65
$ws.Range("N80").Value = -15871.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 4675.36
$ws.Range("I83").Value = 4677.55
$ws.Range("J83").Value = 4666.6
$ws.Range("K83").Value = 42097.95
$ws.Range("L83").Value = 41999.4
$ws.Range("M83").Value = -37417.95
$ws.Range("N83").Value = -51359.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J112").Value = 2000
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 6000
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -8216

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 27500
$ws.Range("J26").Value = 30000
$ws.Range("L26").Value = 30000
$ws.Range("N26").Value = -30560

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H50").Value = 27500
$ws.Range("J50").Value = 30000
$ws.Range("L50").Value = 30000
$ws.Range("N50").Value = -30996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1292.8948
$ws.Range("I102").Value = 1373.5883
$ws.Range("J102").Value = 607
$ws.Range("K102").Value = 1373.5883
$ws.Range("L102").Value = 607
$ws.Range("M102").Value = 248.4117000000001
$ws.Range("N102").Value = -3851

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4772.1665
$ws.Range("I122").Value = 2736.25
$ws.Range("J122").Value = 5790.125
$ws.Range("K122").Value = 8208.75
$ws.Range("L122").Value = 17370.375
$ws.Range("M122").Value = -5758.75
$ws.Range("N122").Value = -22270.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4858
$ws.Range("I126").Value = 1300
$ws.Range("K126").Value = 3900
$ws.Range("M126").Value = -1430

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 4999
$ws.Range("I5").Value = 4999
$ws.Range("K5").Value = 4999
$ws.Range("M5").Value = -4886

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 941.4167
$ws.Range("I22").Value = 829.7
$ws.Range("K22").Value = 829.7
$ws.Range("M22").Value = -534.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 941.4167
$ws.Range("I27").Value = 829.7
$ws.Range("K27").Value = 829.7
$ws.Range("M27").Value = -722.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7709.4
$ws.Range("J40").Value = 8666.333000000001
$ws.Range("L40").Value = 8666.333000000001
$ws.Range("N40").Value = -8938.333000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5467.5454
$ws.Range("I61").Value = 1798.75
$ws.Range("J61").Value = 7564
$ws.Range("K61").Value = 1798.75
$ws.Range("L61").Value = 7564
$ws.Range("M61").Value = -1596.75
$ws.Range("N61").Value = -7968

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 22071.5
$ws.Range("I76").Value = 17999
$ws.Range("K76").Value = 17999
$ws.Range("M76").Value = -17661

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H79").Value = 22071.5
$ws.Range("I79").Value = 17999
$ws.Range("K79").Value = 17999
$ws.Range("M79").Value = -16829

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 5467.5454
$ws.Range("I113").Value = 1798.75
$ws.Range("J113").Value = 7564
$ws.Range("K113").Value = 1798.75
$ws.Range("L113").Value = 7564
$ws.Range("M113").Value = 371.25
$ws.Range("N113").Value = -11904

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4000
$ws.Range("I122").Value = 4000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 12000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -9550
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 6500
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 6500
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 19500
$ws.Range("M126").ClearContents()
